# Update the "Förändrad" (Changed) date column (C) for all data rows
# (rows 2-387) from 2023-09-13 (serial 45182) to 2023-09-15 (serial 45184).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C387").Value = 45184
